# Auto-generated edit script: updates leve profit calculations (currentAveragePrice
# and derived Leve price/profit columns) across all 8 job sheets to reflect refreshed
# market-board data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2970.5518  # "Growing Is Knowing" (was 2785.0908)
$ws.Range("I43").Value = 2986.4092  # "Growing Is Knowing" (was 2787.24)
$ws.Range("J43").Value = 2920.7144  # "Growing Is Knowing" (was 2778.375)
$ws.Range("K43").Value = 2986.4092  # "Growing Is Knowing" (was 2787.24)
$ws.Range("L43").Value = 2920.7144  # "Growing Is Knowing" (was 2778.375)
$ws.Range("M43").Value = -2917.4092  # "Growing Is Knowing" (was -2718.24)
$ws.Range("N43").Value = -3058.7144  # "Growing Is Knowing" (was -2916.375)
$ws.Range("H116").Value = 2072  # "Growing Up" (was 1905.5385)
$ws.Range("I116").Value = 2430.5  # "Growing Up" (was 2224.4)
$ws.Range("J116").Value = 1833  # "Growing Up" (was 1706.25)
$ws.Range("K116").Value = 2430.5  # "Growing Up" (was 2224.4)
$ws.Range("L116").Value = 1833  # "Growing Up" (was 1706.25)
$ws.Range("M116").Value = 1011.5  # "Growing Up" (was 1217.6)
$ws.Range("N116").Value = -8717  # "Growing Up" (was -8590.25)
$ws.Range("H138").Value = 3340.9658  # "All-night Crafting" (was 3410.1785)
$ws.Range("I138").Value = 2761.8667  # "All-night Crafting" (was 2952)
$ws.Range("J138").Value = 3459.959  # "All-night Crafting" (was 3501.8142)
$ws.Range("K138").Value = 8285.6001  # "All-night Crafting" (was 8856)
$ws.Range("L138").Value = 10379.877  # "All-night Crafting" (was 10505.4426)
$ws.Range("M138").Value = -3145.6001  # "All-night Crafting" (was -3716)
$ws.Range("N138").Value = -20659.877  # "All-night Crafting" (was -20785.4426)
$ws.Range("H139").Value = 56985  # "Something Salty and Ceremonial" (was 65880)
$ws.Range("I139").Value = 30000  # "Something Salty and Ceremonial" (was 0)
$ws.Range("J139").Value = 65980  # "Something Salty and Ceremonial" (was 65880)
$ws.Range("K139").Value = 30000  # "Something Salty and Ceremonial" (was 0)
$ws.Range("L139").Value = 65980  # "Something Salty and Ceremonial" (was 65880)
$ws.Range("N139").Value = -76260  # "Something Salty and Ceremonial" (was -76160)
$ws.Range("M139").Value = -24860  # "Something Salty and Ceremonial" (was None)

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 15600  # "Very Slow Array" (was 13008.167)
$ws.Range("J44").Value = 15971.429  # "Very Slow Array" (was 13009.8)
$ws.Range("L44").Value = 15971.429  # "Very Slow Array" (was 13009.8)
$ws.Range("N44").Value = -16947.429  # "Very Slow Array" (was -13985.8)
$ws.Range("H55").Value = 15177.777  # "Employee Retention" (was 12372)
$ws.Range("I55").Value = 9000  # "Employee Retention" (was 0)
$ws.Range("J55").Value = 15950  # "Employee Retention" (was 12372)
$ws.Range("K55").Value = 9000  # "Employee Retention" (was 0)
$ws.Range("L55").Value = 15950  # "Employee Retention" (was 12372)
$ws.Range("N55").Value = -16580  # "Employee Retention" (was -13002)
$ws.Range("M55").Value = -8685  # "Employee Retention" (was None)
$ws.Range("H63").Value = 2998  # "Rivets Run through It" (was 2738)
$ws.Range("J63").Value = 3372.5  # "Rivets Run through It" (was 3047.5)
$ws.Range("L63").Value = 3372.5  # "Rivets Run through It" (was 3047.5)
$ws.Range("N63").Value = -4744.5  # "Rivets Run through It" (was -4419.5)
$ws.Range("H66").Value = 2998  # "A Riveting Revival (L)" (was 2738)
$ws.Range("J66").Value = 3372.5  # "A Riveting Revival (L)" (was 3047.5)
$ws.Range("L66").Value = 16862.5  # "A Riveting Revival (L)" (was 15237.5)
$ws.Range("N66").Value = -23726.5  # "A Riveting Revival (L)" (was -22101.5)
$ws.Range("H122").Value = 3381.6667  # "Haste for High Durium" (was 3682)
$ws.Range("I122").Value = 3004  # "Haste for High Durium" (was 3203)
$ws.Range("J122").Value = 3759.3333  # "Haste for High Durium" (was 4640)
$ws.Range("K122").Value = 9012  # "Haste for High Durium" (was 9609)
$ws.Range("L122").Value = 11277.9999  # "Haste for High Durium" (was 13920)
$ws.Range("M122").Value = -6562  # "Haste for High Durium" (was -7159)
$ws.Range("N122").Value = -16177.9999  # "Haste for High Durium" (was -18820)

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 19833.5  # "Lancers' Creed" (was 19558)
$ws.Range("J35").Value = 19833.5  # "Lancers' Creed" (was 19558)
$ws.Range("L35").Value = 19833.5  # "Lancers' Creed" (was 19558)
$ws.Range("N35").Value = -20453.5  # "Lancers' Creed" (was -20178)
$ws.Range("H80").Value = 1130.5312  # "Unbreaker" (was 1099.3939)
$ws.Range("I80").Value = 681.53845  # "Unbreaker" (was 640.0714)
$ws.Range("J80").Value = 1437.7368  # "Unbreaker" (was 1437.8422)
$ws.Range("K80").Value = 681.53845  # "Unbreaker" (was 640.0714)
$ws.Range("L80").Value = 1437.7368  # "Unbreaker" (was 1437.8422)
$ws.Range("M80").Value = 316.46155  # "Unbreaker" (was 357.9286)
$ws.Range("N80").Value = -3433.7368  # "Unbreaker" (was -3433.8422)
$ws.Range("H83").Value = 1130.5312  # "Attack on Titanium (L)" (was 1099.3939)
$ws.Range("I83").Value = 681.53845  # "Attack on Titanium (L)" (was 640.0714)
$ws.Range("J83").Value = 1437.7368  # "Attack on Titanium (L)" (was 1437.8422)
$ws.Range("K83").Value = 3407.69225  # "Attack on Titanium (L)" (was 3200.357)
$ws.Range("L83").Value = 7188.683999999999  # "Attack on Titanium (L)" (was 7189.211)
$ws.Range("M83").Value = 1584.30775  # "Attack on Titanium (L)" (was 1791.643)
$ws.Range("N83").Value = -17172.684  # "Attack on Titanium (L)" (was -17173.211)
$ws.Range("H86").Value = 226001.2  # "Through Thick and Thin" (was 87830)
$ws.Range("I86").Value = 281501.5  # "Through Thick and Thin" (was 103471.82)
$ws.Range("J86").Value = 4000  # "Through Thick and Thin" (was 1800)
$ws.Range("K86").Value = 281501.5  # "Through Thick and Thin" (was 103471.82)
$ws.Range("L86").Value = 4000  # "Through Thick and Thin" (was 1800)
$ws.Range("M86").Value = -280378.5  # "Through Thick and Thin" (was -102348.82)
$ws.Range("N86").Value = -6246  # "Through Thick and Thin" (was -4046)
$ws.Range("H89").Value = 226001.2  # "Piercing Eyes Deserve Piercing Shafts (L)" (was 87830)
$ws.Range("I89").Value = 281501.5  # "Piercing Eyes Deserve Piercing Shafts (L)" (was 103471.82)
$ws.Range("J89").Value = 4000  # "Piercing Eyes Deserve Piercing Shafts (L)" (was 1800)
$ws.Range("K89").Value = 1407507.5  # "Piercing Eyes Deserve Piercing Shafts (L)" (was 517359.1)
$ws.Range("L89").Value = 20000  # "Piercing Eyes Deserve Piercing Shafts (L)" (was 9000)
$ws.Range("M89").Value = -1401891.5  # "Piercing Eyes Deserve Piercing Shafts (L)" (was -511743.1)
$ws.Range("N89").Value = -31232  # "Piercing Eyes Deserve Piercing Shafts (L)" (was -20232)
$ws.Range("H123").Value = 48000  # "Archon Denied" (was 0)
$ws.Range("J123").Value = 48000  # "Archon Denied" (was 0)
$ws.Range("L123").Value = 48000  # "Archon Denied" (was 0)
$ws.Range("N123").Value = -57800  # "Archon Denied" (was None)

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22824.846  # "Wall Not Found" (was 22786.549)
$ws.Range("I31").Value = 1153.6154  # "Wall Not Found" (was 1136.775)
$ws.Range("J31").Value = 49236.656  # "Wall Not Found" (was 50721.742)
$ws.Range("K31").Value = 1153.6154  # "Wall Not Found" (was 1136.775)
$ws.Range("L31").Value = 49236.656  # "Wall Not Found" (was 50721.742)
$ws.Range("M31").Value = -858.6153999999999  # "Wall Not Found" (was -841.7750000000001)
$ws.Range("N31").Value = -49826.656  # "Wall Not Found" (was -51311.742)
$ws.Range("H34").Value = 22824.846  # "Armoires of the Rich and Famous" (was 22786.549)
$ws.Range("I34").Value = 1153.6154  # "Armoires of the Rich and Famous" (was 1136.775)
$ws.Range("J34").Value = 49236.656  # "Armoires of the Rich and Famous" (was 50721.742)
$ws.Range("K34").Value = 1153.6154  # "Armoires of the Rich and Famous" (was 1136.775)
$ws.Range("L34").Value = 49236.656  # "Armoires of the Rich and Famous" (was 50721.742)
$ws.Range("M34").Value = -951.6153999999999  # "Armoires of the Rich and Famous" (was -934.7750000000001)
$ws.Range("N34").Value = -49640.656  # "Armoires of the Rich and Famous" (was -51125.742)
$ws.Range("H55").Value = 11229.2  # "Ready for a Rematch" (was 13229.2)
$ws.Range("I55").Value = 9036.5  # "Ready for a Rematch" (was 10382)
$ws.Range("J55").Value = 20000  # "Ready for a Rematch" (was 17500)
$ws.Range("K55").Value = 9036.5  # "Ready for a Rematch" (was 10382)
$ws.Range("L55").Value = 20000  # "Ready for a Rematch" (was 17500)
$ws.Range("M55").Value = -8721.5  # "Ready for a Rematch" (was -10067)
$ws.Range("N55").Value = -20630  # "Ready for a Rematch" (was -18130)
$ws.Range("H74").Value = 40016.145  # "License to Heal" (was 40056.668)
$ws.Range("J74").Value = 40016.145  # "License to Heal" (was 40056.668)
$ws.Range("L74").Value = 40016.145  # "License to Heal" (was 40056.668)
$ws.Range("N74").Value = -41764.145  # "License to Heal" (was -41804.668)
$ws.Range("H77").Value = 40016.145  # "Purified Polyrhythm (L)" (was 40056.668)
$ws.Range("J77").Value = 40016.145  # "Purified Polyrhythm (L)" (was 40056.668)
$ws.Range("L77").Value = 120048.435  # "Purified Polyrhythm (L)" (was 120170.004)
$ws.Range("N77").Value = -128784.435  # "Purified Polyrhythm (L)" (was -128906.004)

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 578  # "Flakes for Friends" (was 266.66666)
$ws.Range("J25").Value = 796.6667  # "Flakes for Friends" (was 300)
$ws.Range("L25").Value = 2390.0001  # "Flakes for Friends" (was 900)
$ws.Range("N25").Value = -2728.0001  # "Flakes for Friends" (was -1238)
$ws.Range("H30").Value = 578  # "Picnic Panic" (was 266.66666)
$ws.Range("J30").Value = 796.6667  # "Picnic Panic" (was 300)
$ws.Range("L30").Value = 2390.0001  # "Picnic Panic" (was 900)
$ws.Range("N30").Value = -2594.0001  # "Picnic Panic" (was -1104)
$ws.Range("H131").Value = 756331.0600000001  # "The Mountain Steeped" (was 734434.2)
$ws.Range("J131").Value = 844484.7  # "The Mountain Steeped" (was 817271.9399999999)
$ws.Range("L131").Value = 2533454.1  # "The Mountain Steeped" (was 2451815.82)
$ws.Range("N131").Value = -2543534.1  # "The Mountain Steeped" (was -2461895.82)

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4375.8  # "Awarding Academic Excellence" (was 4415.8)
$ws.Range("I122").Value = 4822.7144  # "Awarding Academic Excellence" (was 4879.857)
$ws.Range("K122").Value = 14468.1432  # "Awarding Academic Excellence" (was 14639.571)
$ws.Range("M122").Value = -12018.1432  # "Awarding Academic Excellence" (was -12189.571)

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2470.2  # "Skin off Their Backs" (was 1533.4445)
$ws.Range("I22").Value = 9800  # "Skin off Their Backs" (was 3533.6667)
$ws.Range("J22").Value = 1655.7778  # "Skin off Their Backs" (was 533.3333)
$ws.Range("K22").Value = 9800  # "Skin off Their Backs" (was 3533.6667)
$ws.Range("L22").Value = 1655.7778  # "Skin off Their Backs" (was 533.3333)
$ws.Range("M22").Value = -9505  # "Skin off Their Backs" (was -3238.6667)
$ws.Range("N22").Value = -2245.7778  # "Skin off Their Backs" (was -1123.3333)
$ws.Range("H27").Value = 2470.2  # "Fire and Hide" (was 1533.4445)
$ws.Range("I27").Value = 9800  # "Fire and Hide" (was 3533.6667)
$ws.Range("J27").Value = 1655.7778  # "Fire and Hide" (was 533.3333)
$ws.Range("K27").Value = 9800  # "Fire and Hide" (was 3533.6667)
$ws.Range("L27").Value = 1655.7778  # "Fire and Hide" (was 533.3333)
$ws.Range("M27").Value = -9693  # "Fire and Hide" (was -3426.6667)
$ws.Range("N27").Value = -1869.7778  # "Fire and Hide" (was -747.3333)
$ws.Range("H46").Value = 5465  # "Supply Side Logic" (was 7120)
$ws.Range("I46").Value = 5546.6665  # "Supply Side Logic" (was 6596)
$ws.Range("J46").Value = 5383.3335  # "Supply Side Logic" (was 7775)
$ws.Range("K46").Value = 5546.6665  # "Supply Side Logic" (was 6596)
$ws.Range("L46").Value = 5383.3335  # "Supply Side Logic" (was 7775)
$ws.Range("M46").Value = -5358.6665  # "Supply Side Logic" (was -6408)
$ws.Range("N46").Value = -5759.3335  # "Supply Side Logic" (was -8151)
$ws.Range("H55").Value = 1415.8857  # "It's Not a Job, It's a Calling" (was 1381.8334)
$ws.Range("J55").Value = 1260.6154  # "It's Not a Job, It's a Calling" (was 1220.963)
$ws.Range("L55").Value = 1260.6154  # "It's Not a Job, It's a Calling" (was 1220.963)
$ws.Range("N55").Value = -1606.6154  # "It's Not a Job, It's a Calling" (was -1566.963)

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1589.7142  # "A Polished Purchase" (was 1605.4286)
$ws.Range("I126").Value = 1688  # "A Polished Purchase" (was 1706.3334)
$ws.Range("K126").Value = 5064  # "A Polished Purchase" (was 5119.0002)
$ws.Range("M126").Value = -2594  # "A Polished Purchase" (was -2649.0002)
$ws.Range("H132").Value = 19173.139  # "Comfy Cabins" (was 21350.21)
$ws.Range("I132").Value = 1894.5834  # "Comfy Cabins" (was 2090.907)
$ws.Range("J132").Value = 102110.2  # "Comfy Cabins" (was 113366.89)
$ws.Range("K132").Value = 5683.7502  # "Comfy Cabins" (was 6272.721)
$ws.Range("L132").Value = 306330.6  # "Comfy Cabins" (was 340100.67)
$ws.Range("M132").Value = -3153.7502  # "Comfy Cabins" (was -3742.721)
$ws.Range("N132").Value = -311390.6  # "Comfy Cabins" (was -345160.67)
